$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P6").Value = 5.15
$ws.Range("Q9").Value = 1.98
$ws.Range("R9").Value = 1.88
$ws.Range("H11").Value = 3.25
$ws.Range("I11").Value = 4.75
$ws.Range("Z11").Value = 11.25
$ws.Range("AH11").Value = 10.5
$ws.Range("H12").Value = 4.5
$ws.Range("W12").Value = 19
$ws.Range("AB12").Value = 51
$ws.Range("AH12").Value = 7
$ws.Range("G13").Value = 2.47
$ws.Range("I13").Value = 2.7
$ws.Range("J13").Value = 3
$ws.Range("L13").Value = 3.2
$ws.Range("Q13").Value = 1.91
$ws.Range("W13").Value = 7.9
$ws.Range("X13").Value = 12
$ws.Range("Y13").Value = 9.5
$ws.Range("Z13").Value = 26
$ws.Range("AA13").Value = 21
$ws.Range("AB13").Value = 30
$ws.Range("AC13").Value = 9.25
$ws.Range("AH13").Value = 8.5
$ws.Range("AI13").Value = 13.5
$ws.Range("AJ13").Value = 10
$ws.Range("AK13").Value = 32
$ws.Range("AL13").Value = 23
$ws.Range("AM13").Value = 32
$ws.Range("AN13").Value = 4.4
$ws.Range("AO13").Value = 13
$ws.Range("AP13").Value = 20
$ws.Range("AQ13").Value = 55
$ws.Range("AR13").Value = 80
$ws.Range("AU13").Value = 6.7
$ws.Range("AX13").Value = 4.65
$ws.Range("AY13").Value = 14
$ws.Range("AZ13").Value = 21
$ws.Range("BA13").Value = 60
$ws.Range("BB13").Value = 90
$ws.Range("G14").Value = 3.15
$ws.Range("H14").Value = 3.05
$ws.Range("I14").Value = 2.25
$ws.Range("J14").Value = 3.65
$ws.Range("K14").Value = 1.98
$ws.Range("L14").Value = 2.92
$ws.Range("M14").Value = 1.03
$ws.Range("N14").Value = 7.7
$ws.Range("O14").Value = 1.29
$ws.Range("P14").Value = 3
$ws.Range("Q14").Value = 1.85
$ws.Range("R14").Value = 1.75
$ws.Range("T14").Value = 2.45
$ws.Range("U14").Value = 1.65
$ws.Range("V14").Value = 2
$ws.Range("W14").Value = 10.25
$ws.Range("X14").Value = 17.5
$ws.Range("Y14").Value = 10.75
$ws.Range("Z14").Value = 45
$ws.Range("AA14").Value = 26
$ws.Range("AB14").Value = 30
$ws.Range("AC14").Value = 9.25
$ws.Range("AD14").Value = 6
$ws.Range("AF14").Value = 55
$ws.Range("AG14").Value = 400
$ws.Range("AH14").Value = 7.7
$ws.Range("AI14").Value = 11.25
$ws.Range("AJ14").Value = 8.75
$ws.Range("AK14").Value = 23
$ws.Range("AL14").Value = 18.5
$ws.Range("AM14").Value = 27
$ws.Range("AN14").Value = 5
$ws.Range("AO14").Value = 17.5
$ws.Range("AP14").Value = 23
$ws.Range("AQ14").Value = 80
$ws.Range("AR14").Value = 110
$ws.Range("AT14").Value = 2.4
$ws.Range("AX14").Value = 4.15
$ws.Range("AY14").Value = 12.5
$ws.Range("AZ14").Value = 21
$ws.Range("BA14").Value = 55
$ws.Range("BB14").Value = 90
$ws.Range("G16").Value = 4.38
$ws.Range("H16").Value = 4.1
$ws.Range("I16").Value = 1.66
$ws.Range("J16").Value = 4.94
$ws.Range("K16").Value = 2.33
$ws.Range("L16").Value = 2.27
$ws.Range("O16").Value = 1.22
$ws.Range("P16").Value = 3.8
$ws.Range("Q16").Value = 1.71
$ws.Range("R16").Value = 2.09
$ws.Range("U16").Value = 1.74
$ws.Range("V16").Value = 2.02
$ws.Range("W16").Value = 11
$ws.Range("X16").Value = 21
$ws.Range("Y16").Value = 11
$ws.Range("Z16").Value = 60
$ws.Range("AA16").Value = 35
$ws.Range("AB16").Value = 35
$ws.Range("AC16").Value = 11
$ws.Range("AD16").Value = 6.2
$ws.Range("AE16").Value = 12
$ws.Range("AF16").Value = 60
$ws.Range("AG16").Value = 101
$ws.Range("AH16").Value = 6.2
$ws.Range("AI16").Value = 6.6
$ws.Range("AJ16").Value = 6.4
$ws.Range("AK16").Value = 10
$ws.Range("AL16").Value = 10
$ws.Range("AM16").Value = 19
$ws.Range("G17").Value = 1.5
$ws.Range("I17").Value = 6.25
$ws.Range("K17").Value = 2.38
$ws.Range("L17").Value = 6
$ws.Range("M17").Value = 1.04
$ws.Range("N17").Value = 12
$ws.Range("Q17").Value = 1.75
$ws.Range("R17").Value = 2.05
$ws.Range("W17").Value = 7
$ws.Range("AC17").Value = 12
$ws.Range("AG17").Value = 301
$ws.Range("AJ17").Value = 19
$ws.Range("AL17").Value = 41
$ws.Range("AM17").Value = 41
$ws.Range("AS17").Value = 126
$ws.Range("AU17").Value = 8.5
$ws.Range("AX17").Value = 7.5
$ws.Range("BC17").Value = 251
$ws.Range("M18").Value = 1.03
$ws.Range("N18").Value = 17
$ws.Range("O18").Value = 1.13
$ws.Range("P18").Value = 6
$ws.Range("Q18").Value = 1.44
$ws.Range("R18").Value = 2.7
$ws.Range("S18").Value = 1.22
$ws.Range("T18").Value = 4
$ws.Range("U18").Value = 2.25
$ws.Range("V18").Value = 1.57
$ws.Range("X18").Value = 6.5
$ws.Range("AA18").Value = 11
$ws.Range("AE18").Value = 29
$ws.Range("AF18").Value = 81
$ws.Range("AG18").Value = 501
$ws.Range("AH18").Value = 41
$ws.Range("AP18").Value = 17
$ws.Range("AS18").Value = 126
$ws.Range("AT18").Value = 4
$ws.Range("AV18").Value = 67
$ws.Range("AY18").Value = 51
$ws.Range("BA18").Value = 401
$ws.Range("BB18").Value = 351
$ws.Range("BC18").Value = 501
$ws.Range("G19").Value = 3.5
$ws.Range("H19").Value = 3.3
$ws.Range("I19").Value = 2.1
$ws.Range("J19").Value = 4.33
$ws.Range("K19").Value = 2
$ws.Range("M19").Value = 1.07
$ws.Range("N19").Value = 9
$ws.Range("U19").Value = 2
$ws.Range("V19").Value = 1.75
$ws.Range("AA19").Value = 34
$ws.Range("AD19").Value = 6.5
$ws.Range("AG19").Value = 401
$ws.Range("AI19").Value = 9
$ws.Range("AP19").Value = 34
$ws.Range("AS19").Value = 301
$ws.Range("AZ19").Value = 26
$ws.Range("G21").Value = 1.52
$ws.Range("I21").Value = 5.4
$ws.Range("J21").Value = 2.07
$ws.Range("K21").Value = 2.25
$ws.Range("L21").Value = 5.5
$ws.Range("N21").Value = 7.9
$ws.Range("S21").Value = 1.37
$ws.Range("V21").Value = 1.85
$ws.Range("Z21").Value = 10.75
$ws.Range("AC21").Value = 7.9
$ws.Range("AE21").Value = 17
$ws.Range("AG21").Value = 600
$ws.Range("AH21").Value = 15.5
$ws.Range("AJ21").Value = 17.5
$ws.Range("AK21").Value = 110
$ws.Range("AL21").Value = 55
$ws.Range("AM21").Value = 55
$ws.Range("AO21").Value = 7.3
$ws.Range("AP21").Value = 17
$ws.Range("AQ21").Value = 23
$ws.Range("AR21").Value = 55
$ws.Range("AS21").Value = 250
$ws.Range("AX21").Value = 7.1
$ws.Range("AY21").Value = 32
$ws.Range("AZ21").Value = 35
$ws.Range("BA21").Value = 200
$ws.Range("BC21").Value = 450
$ws.Range("G25").Value = 1.32
$ws.Range("H25").Value = 4.75
$ws.Range("I25").Value = 8.25
$ws.Range("J25").Value = 1.78
$ws.Range("L25").Value = 7.2
$ws.Range("N25").Value = 14
$ws.Range("O25").Value = 1.22
$ws.Range("P25").Value = 3.45
$ws.Range("Q25").Value = 1.65
$ws.Range("R25").Value = 1.98
$ws.Range("S25").Value = 1.32
$ws.Range("T25").Value = 3.2
$ws.Range("U25").Value = 2.05
$ws.Range("V25").Value = 1.6
$ws.Range("W25").Value = 6.6
$ws.Range("X25").Value = 6
$ws.Range("Y25").Value = 8.75
$ws.Range("Z25").Value = 7.8
$ws.Range("AA25").Value = 11.5
$ws.Range("AB25").Value = 32
$ws.Range("AC25").Value = 12
$ws.Range("AD25").Value = 9.75
$ws.Range("AE25").Value = 24
$ws.Range("AF25").Value = 120
$ws.Range("AH25").Value = 20
$ws.Range("AI25").Value = 55
$ws.Range("AJ25").Value = 27
$ws.Range("AK25").Value = 200
$ws.Range("AL25").Value = 110
$ws.Range("AM25").Value = 100
$ws.Range("AN25").Value = 3.05
$ws.Range("AO25").Value = 5.7
$ws.Range("AP25").Value = 17
$ws.Range("AQ25").Value = 15.5
$ws.Range("AR25").Value = 45
$ws.Range("AS25").Value = 250
$ws.Range("AT25").Value = 2.95
$ws.Range("AU25").Value = 9
$ws.Range("AV25").Value = 90
$ws.Range("AX25").Value = 8.75
$ws.Range("AY25").Value = 50
$ws.Range("AZ25").Value = 50
$ws.Range("BA25").Value = 350
$ws.Range("BB25").Value = 400
